# "Generate Report for Handback"
#
# This localization-status report tracks, per source file / per target
# language, where the file stands in the handoff/handback pipeline. This
# change records that the handback has happened:
#   - every row's Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - each language sheet (zh-cn, de-de) gets its "Latest Target File" and
#     "Latest Handback File" columns (F, G) populated with links to the
#     target markdown file and the handback .xlf translation file
#   - the "Latest Handback DateTime" column (H) is stamped with the time
#     the handback report was generated for that language

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column: every row, every sheet ------------------------------
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- zh-cn: Latest Target File (F) / Latest Handback File (G) -----------
$zhcnRepoRoot = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c574756a97bd4ec7b4403cb6b9cbca16e48e5fb3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/hb/"
$zhcnXlfName  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/109883b34df792143840eef21bd76471f7cdd7ab/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), ($zhcnRepoRoot + $zhcnXlfName), "", "", $zhcnXlfName)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/109883b34df792143840eef21bd76471f7cdd7ab/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), ($zhcnRepoRoot + $zhcnXlfName), "", "", $zhcnXlfName)

$zhcn.Range("F2:G2").Font.Underline = $true
$zhcn.Range("F2:G2").Font.Color = 15570276
$zhcn.Range("F3:G3").Font.Underline = $true
$zhcn.Range("F3:G3").Font.Color = 15570276

# Latest Handback DateTime (H) for zh-cn
$zhcn.Range("H2").Value = "2016-03-14 09:12:48"
$zhcn.Range("H3").Value = "2016-03-14 09:12:48"

# --- de-de: Latest Target File (F) / Latest Handback File (G) -----------
$dedeRepoRoot = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c2ba04b5d7d92c34aa3c375c298c8d015a6e5525/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/hb/"
$dedeXlfName  = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/109883b34df792143840eef21bd76471f7cdd7ab/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G2"), ($dedeRepoRoot + $dedeXlfName), "", "", $dedeXlfName)
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/109883b34df792143840eef21bd76471f7cdd7ab/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("G3"), ($dedeRepoRoot + $dedeXlfName), "", "", $dedeXlfName)

$dede.Range("F2:G2").Font.Underline = $true
$dede.Range("F2:G2").Font.Color = 15570276
$dede.Range("F3:G3").Font.Underline = $true
$dede.Range("F3:G3").Font.Color = 15570276

# Latest Handback DateTime (H) for de-de
$dede.Range("H2").Value = "2016-03-14 09:13:01"
$dede.Range("H3").Value = "2016-03-14 09:13:01"
